$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="39.767.52"},
    @{Cell="E2"; Value="  -4.66%  "},
    @{Cell="D3"; Value="2.323.09"},
    @{Cell="E3"; Value="  -6.19%  "},
    @{Cell="E4"; Value="  -0.13%  "},
    @{Cell="D5"; Value="306.05"},
    @{Cell="E5"; Value="  -4.45%  "},
    @{Cell="D6"; Value="83.69"},
    @{Cell="E6"; Value="  -9.36%  "},
    @{Cell="D7"; Value="0.525"},
    @{Cell="E7"; Value="  -4.62%  "},
    @{Cell="E8"; Value="  +0.02%  "},
    @{Cell="D9"; Value="0.482"},
    @{Cell="E9"; Value="  -5.73%  "},
    @{Cell="D10"; Value="0.0817"},
    @{Cell="E10"; Value="  -5.35%  "},
    @{Cell="D11"; Value="29.82"},
    @{Cell="E11"; Value="  -9.77%  "},
    @{Cell="E12"; Value="  -0.31%  "},
    @{Cell="D13"; Value="2.682.05"},
    @{Cell="E13"; Value="  -6.17%  "},
    @{Cell="D14"; Value="6.37"},
    @{Cell="E14"; Value="  -7.76%  "},
    @{Cell="D15"; Value="14.70"},
    @{Cell="E15"; Value="  -5.09%  "},
    @{Cell="D16"; Value="2.333.58"},
    @{Cell="E16"; Value="  -5.88%  "},
    @{Cell="D17"; Value="0.743"},
    @{Cell="E17"; Value="  -6.49%  "},
    @{Cell="D18"; Value="39.708.16"},
    @{Cell="E18"; Value="  -4.61%  "},
    @{Cell="D19"; Value="0.0₃0899"},
    @{Cell="E19"; Value="  -4.67%  "},
    @{Cell="D20"; Value="6.04"},
    @{Cell="E20"; Value="  -6.30%  "},
    @{Cell="D21"; Value="67.45"},
    @{Cell="E21"; Value="  -4.55%  "},
    @{Cell="D22"; Value="10.55"},
    @{Cell="E22"; Value="  -6.34%  "},
    @{Cell="D23"; Value="233.86"},
    @{Cell="E23"; Value="  -2.50%  "},
    @{Cell="D24"; Value="2.53"},
    @{Cell="E24"; Value="  -8.34%  "},
    @{Cell="E25"; Value="  +0.32%  "},
    @{Cell="D26"; Value="1.79"},
    @{Cell="E26"; Value="  -8.19%  "},
    @{Cell="D27"; Value="23.22"},
    @{Cell="E27"; Value="  -7.19%  "},
    @{Cell="D28"; Value="2.21"},
    @{Cell="E28"; Value="  -1.64%  "},
    @{Cell="D29"; Value="9.14"},
    @{Cell="E29"; Value="  -6.30%  "},
    @{Cell="D30"; Value="34.22"},
    @{Cell="E30"; Value="  -6.86%  "},
    @{Cell="D31"; Value="150.02"},
    @{Cell="E31"; Value="  -4.84%  "},
    @{Cell="D32"; Value="0.999"},
    @{Cell="E32"; Value="  -0.25%  "},
    @{Cell="D33"; Value="5.08"},
    @{Cell="E33"; Value="  -6.62%  "},
    @{Cell="B34"; Value="Hedera"},
    @{Cell="C34"; Value="https://coinranking.com/coin/jad286TjB+hedera-hbar"},
    @{Cell="D34"; Value="0.0720"},
    @{Cell="E34"; Value="  -5.85%  "},
    @{Cell="B35"; Value="WEMIXToken"},
    @{Cell="C35"; Value="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"},
    @{Cell="D35"; Value="2.42"},
    @{Cell="E35"; Value="  -5.78%  "},
    @{Cell="D36"; Value="0.113"},
    @{Cell="E36"; Value="  -3.05%  "},
    @{Cell="D37"; Value="2.74"},
    @{Cell="E37"; Value="  -4.97%  "},
    @{Cell="D38"; Value="0.0986"},
    @{Cell="E38"; Value="  -4.90%  "},
    @{Cell="D39"; Value="15.64"},
    @{Cell="E39"; Value="  -9.09%  "},
    @{Cell="D40"; Value="1.69"},
    @{Cell="E40"; Value="  -8.60%  "},
    @{Cell="D41"; Value="3.77"},
    @{Cell="E41"; Value="  -6.61%  "},
    @{Cell="E42"; Value="  -5.55%  "},
    @{Cell="D43"; Value="1.933.56"},
    @{Cell="E43"; Value="  -3.06%  "},
    @{Cell="D44"; Value="0.0263"},
    @{Cell="E44"; Value="  -7.34%  "},
    @{Cell="D45"; Value="17.42"},
    @{Cell="E45"; Value="  -7.21%  "},
    @{Cell="D46"; Value="9.20"},
    @{Cell="E46"; Value="  -3.08%  "},
    @{Cell="D47"; Value="2.63"},
    @{Cell="E47"; Value="  -11.77%  "},
    @{Cell="D48"; Value="2.546.55"},
    @{Cell="E48"; Value="  -7.50%  "},
    @{Cell="D49"; Value="91.74"},
    @{Cell="E49"; Value="  -5.97%  "},
    @{Cell="D50"; Value="70.40"},
    @{Cell="E50"; Value="  -7.38%  "},
    @{Cell="D51"; Value="63.11"},
    @{Cell="E51"; Value="  -6.52%  "}
)

foreach ($item in $changes) {
    $ws.Range($item.Cell).NumberFormat = "@"
    $ws.Range($item.Cell).Value = $item.Value
    $ws.Range($item.Cell).ClearFormats()
}
